$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Set D (Price) and E (Volume 1h) columns as plain text to avoid Excel
# auto-converting numeric-looking strings (e.g. "1.001") into numbers.

$dCells = @("D2", "D3", "D4", "D5", "D6", "D7", "D8", "D9", "D10", "D11", "D12", "D13", "D14", "D15", "D16", "D17", "D18", "D20", "D21", "D22", "D23", "D24", "D25", "D26", "D27", "D28", "D29", "D30", "D31", "D32", "D33", "D34", "D35", "D36", "D37", "D38", "D39", "D40", "D41", "D42", "D43", "D44", "D45", "D46", "D47", "D48", "D49", "D50", "D51")
foreach ($addr in $dCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D2").Value = "27.316.32"
$ws.Range("E2").Value = "  +0.80%  "
$ws.Range("D3").Value = "1.903.91"
$ws.Range("E3").Value = "  +0.63%  "
$ws.Range("D4").Value = "1.001"
$ws.Range("E4").Value = "  +0.03%  "
$ws.Range("D5").Value = "308.24"
$ws.Range("E5").Value = "  +0.48%  "
$ws.Range("D6").Value = "1.000"
$ws.Range("E6").Value = "  +0.00%  "
$ws.Range("D7").Value = "0.5220"
$ws.Range("E7").Value = "  +0.57%  "
$ws.Range("D8").Value = "0.3782"
$ws.Range("E8").Value = "  +0.57%  "
$ws.Range("D9").Value = "0.07299"
$ws.Range("E9").Value = "  +1.05%  "
$ws.Range("D10").Value = "21.33"
$ws.Range("E10").Value = "  +0.79%  "
$ws.Range("D11").Value = "0.9039"
$ws.Range("E11").Value = "  +0.14%  "
$ws.Range("D12").Value = "0.08254"
$ws.Range("E12").Value = "  +7.73%  "
$ws.Range("D13").Value = "97.12"
$ws.Range("E13").Value = "  +2.81%  "
$ws.Range("D14").Value = "1.903.02"
$ws.Range("E14").Value = "  +0.67%  "
$ws.Range("D15").Value = "5.302"
$ws.Range("E15").Value = "  +1.14%  "
$ws.Range("D16").Value = "1.001"
$ws.Range("E16").Value = "  +0.02%  "
$ws.Range("D17").Value = "0.000008652"
$ws.Range("E17").Value = "  +1.67%  "
$ws.Range("D18").Value = "14.60"
$ws.Range("E18").Value = "  +1.09%  "
$ws.Range("E19").Value = "  +0.03%  "
$ws.Range("D20").Value = "27.336.48"
$ws.Range("E20").Value = "  +0.70%  "
$ws.Range("D21").Value = "5.102"
$ws.Range("E21").Value = "  +0.72%  "
$ws.Range("D22").Value = "10.69"
$ws.Range("E22").Value = "  +0.85%  "
$ws.Range("D23").Value = "6.447"
$ws.Range("E23").Value = "  +0.94%  "
$ws.Range("D24").Value = "2.310"
$ws.Range("E24").Value = "  +0.02%  "
$ws.Range("D25").Value = "147.61"
$ws.Range("E25").Value = "  +1.31%  "
$ws.Range("D26").Value = "18.27"
$ws.Range("E26").Value = "  +1.11%  "
$ws.Range("D27").Value = "1.748"
$ws.Range("E27").Value = "  +1.18%  "
$ws.Range("D28").Value = "115.51"
$ws.Range("E28").Value = "  +0.90%  "
$ws.Range("D29").Value = "4.856"
$ws.Range("D30").Value = "4.933"
$ws.Range("E30").Value = "  -0.39%  "
$ws.Range("D31").Value = "0.09251"
$ws.Range("E31").Value = "  +0.41%  "
$ws.Range("D32").Value = "0.05076"
$ws.Range("E32").Value = "  +0.39%  "
$ws.Range("D33").Value = "0.8003"
$ws.Range("E33").Value = "  +3.46%  "
$ws.Range("D34").Value = "1.238"
$ws.Range("E34").Value = "  -0.81%  "
$ws.Range("D35").Value = "3.443"
$ws.Range("E35").Value = "  +4.76%  "
$ws.Range("D36").Value = "2.948"
$ws.Range("E36").Value = "  -1.01%  "
$ws.Range("D37").Value = "2.601"
$ws.Range("E37").Value = "  -0.06%  "
$ws.Range("D38").Value = "0.5728"
$ws.Range("E38").Value = "  +0.86%  "
$ws.Range("D39").Value = "0.02007"
$ws.Range("E39").Value = "  +0.73%  "
$ws.Range("D40").Value = "1.078"
$ws.Range("E40").Value = "  +0.43%  "
$ws.Range("D41").Value = "9.034"
$ws.Range("E41").Value = "  -0.21%  "
$ws.Range("D42").Value = "6.591"
$ws.Range("E42").Value = "  -0.73%  "
$ws.Range("D43").Value = "116.27"
$ws.Range("E43").Value = "  -2.77%  "
$ws.Range("D44").Value = "0.1522"
$ws.Range("E44").Value = "  +0.61%  "
$ws.Range("D45").Value = "0.4890"
$ws.Range("E45").Value = "  +0.59%  "
$ws.Range("D46").Value = "1.000"
$ws.Range("E46").Value = "  -0.04%  "
$ws.Range("D47").Value = "10.10"
$ws.Range("E47").Value = "  -0.45%  "
$ws.Range("D48").Value = "1.630"
$ws.Range("E48").Value = "  +1.98%  "
$ws.Range("D49").Value = "38.07"
$ws.Range("E49").Value = "  +1.00%  "
$ws.Range("D50").Value = "63.96"
$ws.Range("E50").Value = "  -0.40%  "
$ws.Range("D51").Value = "0.05947"
$ws.Range("E51").Value = "  +0.50%  "

foreach ($addr in $dCells) {
    $ws.Range($addr).Style = "Normal"
}
